$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:C1) keeps left/right/top thin borders but loses its bottom border
$ws.Range("A1:C1").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

# Data rows 3-29 (A:C) lose their border entirely
$ws.Range("A3:C29").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

# Update the active selection
$ws.Range("C8").Select() | Out-Null
